$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D column (UOM) values from NOS to KG first (matches shared-string insertion order)
$ws.Range("D2").Value = "KG"
$ws.Range("D3").Value = "KG"
$ws.Range("D4").Value = "KG"

# Update B column (FACTORY_CODE) values for data rows
$ws.Range("B2").Value = "IND_TNA_Unit_1"
$ws.Range("B3").Value = "IN_MA_TH_Unit1"
$ws.Range("B4").Value = "__DA_Unit2"

# Update selection to match target workbook view
[void]$ws.Range("B10").Select()
